$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-22 20:58:52"

for ($row = 2; $row -le 73; $row++) {
    $ws.Range("O$row").Value = $newTimestamp
}
